# Brisbane testing sketch added to testing branch:
# Adds a new "os_boundary" column (C) to the open_space_defs sheet, with
# a handful of new lookup values (recreation_ground/village_green/
# cemetary/boundary) layered into the existing os_landuse /
# os_add_as_tags lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("open_space_defs")

# --- Insert a new column before the existing column C (specific_inclusion_criteria) ---
$ws.Range("C1").EntireColumn.Insert()
# Match the width of the neighbouring os_landuse column so it renders the same way.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Phase 1: author the brand new lookup values, in the order they were
#     first typed, so new entries line up with how the sheet was edited ---
$ws.Range("B17").Value = "recreation_ground"
$ws.Range("B20").Value = "village_green"
$ws.Range("B3").Value  = "cemetary"

$ws.Range("C1").Value  = "os_boundary"
$ws.Range("C4").Value  = "nature_reserve"
$ws.Range("C2").Value  = "protected_area"
$ws.Range("C3").Value  = "national_park"
$ws.Range("C6").Value  = "water_protection_area"
$ws.Range("J6").Value  = "boundary"
$ws.Range("C7").Value  = "state_forest"
$ws.Range("C8").Value  = "state_park"
$ws.Range("C9").Value  = "regional_park"
$ws.Range("C11").Value = "county_park"

# --- Phase 2: re-lay the rest of the three lists underneath/around the
#     newly inserted entries (these all reuse already-known values) ---

# Column A (possible_os_tags): add "boundary" entry at the end.
$ws.Range("A13").Value = "boundary"

# Column B (os_landuse): the remainder of the list, shifted down to make
# room for cemetary / recreation_ground / village_green above.
$ws.Range("B4").Value  = "conservation"
$ws.Range("B5").Value  = "field"
$ws.Range("B6").Value  = "forest"
$ws.Range("B7").Value  = "garden"
$ws.Range("B8").Value  = "grass"
$ws.Range("B9").Value  = "green"
$ws.Range("B10").Value = "leisure"
$ws.Range("B11").Value = "meadow"
$ws.Range("B12").Value = "orchard"
$ws.Range("B13").Value = "park"
$ws.Range("B14").Value = "pitch"
$ws.Range("B15").Value = "pond"
$ws.Range("B16").Value = "recreation ground"
$ws.Range("B18").Value = "sport"
$ws.Range("B19").Value = "trees"
$ws.Range("B21").Value = "village green"
$ws.Range("B22").Value = "water"
$ws.Range("B23").Value = "winter_sports"
$ws.Range("B24").Value = "wood"

# Column C (os_boundary): fill in the remaining reused value.
$ws.Range("C5").Value  = "forest"
$ws.Range("C10").Value = "park"

# Column J (os_add_as_tags, previously I): the remainder of the list,
# shifted down to make room for the new "boundary" tag.
$ws.Range("J7").Value  = "landuse"
$ws.Range("J8").Value  = "leisure"
$ws.Range("J9").Value  = "natural"
$ws.Range("J10").Value = "sport"
$ws.Range("J11").Value = "waterway"
$ws.Range("J12").Value = "wood"
$ws.Range("J13").Value = "in_school"
$ws.Range("J14").Value = "is_school"
$ws.Range("J15").Value = "water_feature"
$ws.Range("J16").Value = "medial_axis_length"
$ws.Range("J17").Value = "num_symdiff_convhull_geoms"
$ws.Range("J18").Value = "roundness"
$ws.Range("J19").Value = "linear_feature"
$ws.Range("J20").Value = "acceptable_linear_feature"

# --- Restore view state: active selection on J14 ---
$ws.Activate()
$ws.Range("J14").Select()
